$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Correct Predictions (B) and Incorrect Predictions (C) per row
$ws.Range("B2").Value = 3092
$ws.Range("C2").Value = 3071

$ws.Range("B3").Value = 1791
$ws.Range("C3").Value = 960

$ws.Range("B4").Value = 3398
$ws.Range("C4").Value = 1994

$ws.Range("B5").Value = 807
$ws.Range("C5").Value = 1582

$ws.Range("B6").Value = 1084
$ws.Range("C6").Value = 957

$ws.Range("B7").Value = 1891
$ws.Range("C7").Value = 2785

# Update Cross-Entropy Loss (D) for all rows to the new constant value
$ws.Range("D2:D7").Value = 1.546168446540833
